$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Price" column (D) values look numeric (e.g. "299.44") but must stay as literal
# text, matching the source data feed formatting (some contain multiple dots,
# e.g. "23.433.23"). Force text storage via NumberFormat "@" before assigning,
# then reset the style to Normal so no stray number-format style lingers on the cell.
$priceUpdates = [ordered]@{
    'D2' = '23.433.23'
    'D6' = '299.44'
    'D8' = '0.3556'
    'D9' = '49.78'
    'D10' = '0.08099'
    'D14' = '6.392'
    'D15' = '7.360'
    'D16' = '0.00001196'
    'D17' = '1.651.38'
    'D18' = '97.39'
    'D19' = '0.06942'
    'D20' = '6.766'
    'D24' = '23.459.42'
    'D25' = '2.500'
    'D26' = '2.918'
    'D27' = '20.90'
    'D28' = '152.77'
    'D29' = '5.205'
    'D30' = '132.83'
    'D31' = '1.825.55'
    'D32' = '6.927'
    'D33' = '2.096'
    'D34' = '11.73'
    'D35' = '1.009'
    'D36' = '0.02721'
    'D37' = '0.08730'
    'D38' = '0.2432'
    'D39' = '5.924'
    'D40' = '13.04'
    'D41' = '0.06768'
    'D42' = '0.6877'
    'D43' = '1.303'
    'D44' = '15.57'
    'D46' = '0.6392'
    'D47' = '2.250'
    'D49' = '0.07720'
    'D50' = '127.24'
    'D51' = '1.148'
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = "Normal"
}

# "Volume(1h)" column (E) values are space-padded percentages (e.g. "  -1.23%  ")
# which Excel always treats as text, so a plain assignment is sufficient.
$volumeUpdates = [ordered]@{
    'E2' = '  -1.23%  '
    'E3' = '  -0.30%  '
    'E4' = '  +0.21%  '
    'E5' = '  +0.19%  '
    'E6' = '  -1.51%  '
    'E7' = '  -0.87%  '
    'E8' = '  -1.43%  '
    'E9' = '  -3.11%  '
    'E10' = '  -1.67%  '
    'E11' = '  -2.31%  '
    'E12' = '  +0.12%  '
    'E13' = '  -2.66%  '
    'E14' = '  -2.33%  '
    'E15' = '  -0.65%  '
    'E16' = '  -3.08%  '
    'E17' = '  -0.19%  '
    'E18' = '  +0.08%  '
    'E19' = '  -0.32%  '
    'E20' = '  -0.21%  '
    'E21' = '  -2.29%  '
    'E22' = '  +0.09%  '
    'E23' = '  -1.71%  '
    'E24' = '  -1.12%  '
    'E25' = '  -2.21%  '
    'E26' = '  -5.06%  '
    'E27' = '  -1.98%  '
    'E28' = '  +0.69%  '
    'E29' = '  -0.91%  '
    'E30' = '  -1.66%  '
    'E31' = '  -0.50%  '
    'E32' = '  +0.95%  '
    'E33' = '  -0.64%  '
    'E34' = '  -1.30%  '
    'E35' = '  -7.75%  '
    'E36' = '  -3.84%  '
    'E37' = '  -1.26%  '
    'E38' = '  -3.53%  '
    'E39' = '  -2.81%  '
    'E40' = '  +1.56%  '
    'E42' = '  -2.70%  '
    'E43' = '  -2.60%  '
    'E44' = '  -3.02%  '
    'E45' = '  +0.12%  '
    'E46' = '  -2.17%  '
    'E47' = '  -3.89%  '
    'E48' = '  -1.71%  '
    'E49' = '  -3.33%  '
    'E50' = '  -0.71%  '
    'E51' = '  -3.49%  '
}

foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}
